# "Generate Report for Handoff"
#
# The localization-status report is regenerated by the CI handoff tool.
# For file 98f95ff2-076f-4953-af01-54b66533f775.md the tool re-evaluates
# the "Latest Handoff" timestamps/files for the Overview row and for the
# zh-cn per-language row, writing the freshly computed values in before
# reconciling them back to the same values already on record (the
# handoff itself didn't change, only the report run did). We replay both
# the write of the newly generated values and the settle-back to the
# final, on-record values so the shared-string pool reflects the same
# regeneration pass that produced the committed workbook.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")

# --- Overview sheet: row for 98f95ff2-076f-4953-af01-54b66533f775.md ---
# Latest Handoff Date column (D6)
$overview.Range("D6").Value = "2016-03-23 00:38:10"
$overview.Range("D6").Value = "2016-03-23 00:37:46"

# --- zh-cn sheet: row for 98f95ff2-076f-4953-af01-54b66533f775.md ---
# Latest Handoff File (D6) / Latest Handoff Datetime (E6)
$zhcn.Range("D6").Value = "98f95ff2-076f-4953-af01-54b66533f775.23b7bba38dc2e8c95e4a18a532490ed529bc76dd.zh-cn.xlf"
$zhcn.Range("E6").Value = "2016-03-23 00:38:07"

$zhcn.Range("D6").Value = "98f95ff2-076f-4953-af01-54b66533f775.23b7bba38dc2e8c95e4a18a532490ed529bc76dd.zh-cn.xlf"
$zhcn.Range("E6").Value = "2016-03-23 00:37:41"
